$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 191; this shifts the existing rows
# 191-230 down to 192-231 (dimension grows from A1:R230 to A1:R231).
$ws.Rows("191:191").Insert()

# Populate the newly inserted row 191 with the new data record.
$ws.Cells.Item(191, 1).Value = 8
$ws.Cells.Item(191, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(191, 3).Value = "Coquimbo"
$ws.Cells.Item(191, 4).Value = 44711
$ws.Cells.Item(191, 5).Value = 4
$ws.Cells.Item(191, 6).Value = 100112031
$ws.Cells.Item(191, 7).Value = "Poroto verde"
$ws.Cells.Item(191, 8).Value = "Magnum"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 440
$ws.Cells.Item(191, 11).Value = 24000
$ws.Cells.Item(191, 12).Value = 25000
$ws.Cells.Item(191, 13).Value = 24500
$ws.Cells.Item(191, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(191, 15).Value = "Perú"
$ws.Cells.Item(191, 16).Value = 980
$ws.Cells.Item(191, 17).Value = 25
$ws.Cells.Item(191, 18).Value = "Hortaliza"
